# Appointment ID 8 (row 9) had its Status temporarily marked "Confirmed"
# while trying to resolve the doctor's appointment, but that could not be
# finalized, so the Status is reverted back to "Pending".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G9").Value = "Pending"
